$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 12.6326
$ws.Range("C9").Value = -11.9866
$ws.Range("E11").Value = 13.4865
$ws.Range("C18").Value = -14.3998
$ws.Range("C20").Value = -13.74389999999999
